# Weekly update: a new week of "Jengibre" price data was reported for
# Terminal La Palmera de La Serena. The new observation is inserted as a
# new row 102 (most-recent-first ordering), pushing all the existing data
# rows (previously 102:142) down by one row to 103:143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 102. This shifts rows 102:142 down
# to become rows 103:143 (and grows the sheet dimension to A1:R143).
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly record.
$ws.Range("A102").Value = 8
$ws.Range("B102").Value = "Terminal La Palmera de La Serena"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 45097
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100114007
$ws.Range("G102").Value = "Jengibre"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 240
$ws.Range("K102").Value = 17000
$ws.Range("L102").Value = 18000
$ws.Range("M102").Value = 17500
$ws.Range("N102").Value = "`$/caja 13 kilos"
$ws.Range("O102").Value = "Perú"
$ws.Range("P102").Value = 1346
$ws.Range("Q102").Value = 13
$ws.Range("R102").Value = "Hortaliza"
